$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values to be inserted into column B for rows 2..20 (quarter-ahead naive
# forecast error, previously missing a leading column causing every series to
# be off-by-one quarter).
$newB = @{
    2  = -0.58442257821662
    3  = 1.454533757567239
    4  = -1.777394389465022
    5  = -1.722070219091221
    6  = 0.3648791949059138
    7  = -0.2352699264540507
    8  = -0.05148746350304451
    9  = -0.1333319740152609
    10 = 1.614150253737389
    11 = 0.5701030647716323
    12 = 0.2202779152847414
    13 = 0.5040960054549828
    14 = 0.420735823599318
    15 = -0.1252583916527783
    16 = 0.08824118641116785
    17 = -0.1133200159455487
    18 = 0.1743923273248104
    19 = -0.4559694969238889
    20 = 0.1808172637304477
}

# Shift the existing data in columns B..J right into C..K (column K's old
# value is pushed off the end and discarded), for every data row. Walk the
# columns from the right (K) towards the left (C) so each destination is
# written before its own value is needed as a source.
for ($r = 2; $r -le 20; $r++) {
    for ($c = 11; $c -ge 3; $c--) {
        $srcVal = $ws.Cells.Item($r, $c - 1).Value2
        $ws.Cells.Item($r, $c).Value2 = $srcVal
    }
    $ws.Cells.Item($r, 2).Value2 = $newB[$r]
}
